$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6060606060606061
$ws.Range("P2").Value = 0.1515151515151515
$ws.Range("S2").Value = 0.06060606060606061

$ws.Range("P3").Value = 0.8421052631578947
$ws.Range("S3").Value = 0.1578947368421053

$ws.Range("B6").Value = 0.1176470588235294
$ws.Range("D6").Value = 0.05882352941176471
$ws.Range("F6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.1764705882352941
$ws.Range("Q6").Value = 0.1764705882352941
$ws.Range("R6").Value = 0.05882352941176471
$ws.Range("S6").Value = 0.3529411764705883

$ws.Range("B7").Value = 0.1071428571428571
$ws.Range("F7").Value = 0.1785714285714286
$ws.Range("J7").Value = 0.03571428571428571
$ws.Range("O7").Value = 0.07142857142857142
$ws.Range("Q7").Value = 0.03571428571428571
$ws.Range("R7").Value = 0.03571428571428571
$ws.Range("S7").Value = 0.5357142857142857

$ws.Range("B8").Value = 0.1730769230769231
$ws.Range("F8").Value = 0.0576923076923077
$ws.Range("J8").Value = 0.1153846153846154
$ws.Range("O8").Value = 0.01923076923076923
$ws.Range("Q8").Value = 0.09615384615384616
$ws.Range("R8").Value = 0.03846153846153846
$ws.Range("S8").Value = 0.5

$ws.Range("B9").Value = 0.1666666666666667
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("S9").Value = 0.5

$ws.Range("B10").Value = 0.1071428571428571
$ws.Range("F10").Value = 0.05952380952380952
$ws.Range("J10").Value = 0.119047619047619
$ws.Range("Q10").Value = 0.1547619047619048
$ws.Range("R10").Value = 0.04761904761904762
$ws.Range("S10").Value = 0.5119047619047619

$ws.Range("G11").Value = 0.1612903225806452
$ws.Range("J11").Value = 0.06451612903225806
$ws.Range("K11").Value = 0.1612903225806452
$ws.Range("L11").Value = 0.6129032258064516

$ws.Range("G12").Value = 0.8571428571428571
$ws.Range("J12").Value = 0.04761904761904762
$ws.Range("L12").Value = 0.09523809523809523

$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.1428571428571428
$ws.Range("S13").Value = 0.1428571428571428

$ws.Range("H15").Value = 0.08333333333333333
$ws.Range("J15").Value = 0.25
$ws.Range("K15").Value = 0.08333333333333333
$ws.Range("S15").Value = 0.5833333333333334

$ws.Range("F16").Value = 0.04545454545454546
$ws.Range("H16").Value = 0.2727272727272727
$ws.Range("I16").Value = 0.04545454545454546
$ws.Range("J16").Value = 0.3636363636363636
$ws.Range("K16").Value = 0.04545454545454546
$ws.Range("O16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.1363636363636364

$ws.Range("H17").Value = 0.2307692307692308
$ws.Range("I17").Value = 0.1538461538461539
$ws.Range("J17").Value = 0.3461538461538461
$ws.Range("K17").Value = 0.1153846153846154
$ws.Range("M17").Value = 0.03846153846153846
$ws.Range("O17").Value = 0.03846153846153846
$ws.Range("S17").Value = 0.07692307692307693

$ws.Range("I18").Value = 0.25
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.125
$ws.Range("S18").Value = 0.125

$ws.Range("H19").Value = 0.2727272727272727
$ws.Range("I19").Value = 0.07692307692307693
$ws.Range("J19").Value = 0.2517482517482518
$ws.Range("K19").Value = 0.1398601398601399
$ws.Range("M19").Value = 0.04195804195804196
$ws.Range("O19").Value = 0.02797202797202797
$ws.Range("S19").Value = 0.1888111888111888
